# FX module - front panel adjustments
# Update the two RV1 "X" readings (column C) for the SW1 / SW2 rows;
# the dependent shared-formula cells (E5:E6, I5:I6) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = 61.142380000000003
$ws.Range("C6").Value = 61.143000000000001

# Move the front-panel selection to C6, matching where the edit was made.
$ws.Range("C6").Select()
